$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell L1 -> _MasterItemLabel (new shared string, index 34)
$ws.Range("L1").Value = "_MasterItemLabel"

# New data cells in column L for rows 3 and 5
$ws.Range("L3").Value = "GetSelectedCount(Alpha)"
$ws.Range("L5").Value = "Dim Name Here!"

# Update the active selection from L4 to L6
$ws.Range("L6").Select()
